# Apply the latest cryptos-list refresh (prices + 1h volume deltas,
# plus a handful of re-ranked coins whose rows traded places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain string/text updates -- none of these parse as a plain number,
# so Excel keeps them as literal text on assignment.
$ws.Range("D2").Value = '66.328.71'
$ws.Range("E2").Value = '  +1.87%  '
$ws.Range("D3").Value = '3.223.19'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +4.28%  '
$ws.Range("E6").Value = '  +2.55%  '
$ws.Range("D8").Value = '3.223.85'
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("E9").Value = '  +1.24%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("D15").Value = '3.749.16'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("E16").Value = '  +4.10%  '
$ws.Range("D17").Value = '66.382.84'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").Value = '3.231.99'
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  +5.22%  '
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E27").Value = '  +3.25%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("E28").Value = '  +3.57%  '
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("E30").Value = '  +10.75%  '
$ws.Range("E31").Value = '  +3.63%  '
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("E37").Value = '  +2.10%  '
$ws.Range("E38").Value = '  +2.28%  '
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("E40").Value = '  -3.91%  '
$ws.Range("E41").Value = '  +2.48%  '
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("E42").Value = '  +3.10%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E44").Value = '  +4.88%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.977.15'
$ws.Range("E45").Value = '  -2.96%  '
$ws.Range("E46").Value = '  +5.65%  '
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("E49").Value = '  +3.73%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("E51").Value = '  +5.38%  '

# The "Price" column stores numeric-looking values (e.g. "604.31",
# "0.120") as TEXT in the workbook, so trailing zeros and plain
# decimals survive. Assigning such a string straight to .Value would
# make Excel auto-convert it to a real number (losing formatting like
# "155.30" -> 155.3). To keep these as text without perturbing the
# cell style, stage each value in a scratch cell that is explicitly
# text-formatted, then copy/paste-special (values only) into place --
# PasteSpecial keeps the destination cells own (default) style.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = '604.31'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = '155.30'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = '0.538'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Value = '0.163'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Value = '6.19'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Value = '0.512'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$scratch.Value = '0.0000275'
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Value = '7.49'
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Value = '514.90'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$scratch.Value = '15.70'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Value = '0.741'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Value = '85.64'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Value = '9.34'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Value = '3.05'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Value = '2.26'
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Value = '28.38'
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Value = '55.53'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Value = '0.0922'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Value = '489.03'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Value = '0.0425'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$scratch.Value = '8.89'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Value = '0.298'
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Value = '0.120'
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$scratch.Value = '2.51'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$scratch.Value = '29.29'
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Value = '2.36'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Value = '0.117'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Value = '34.24'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
